$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended below the existing table (row 55).
# The date column stores a literal text string (e.g. "2025/10/03"), matching
# the existing rows above it, rather than Excel's auto-recognized date
# serial number. Force the cell to Text format before assigning the value so
# Excel keeps it as a string, then clear the temporary formatting so the
# cell is left with the sheet's default (unstyled) appearance, same as its
# neighbours.
$ws.Range("A55").NumberFormat = "@"
$ws.Range("A55").Value = "2025/10/03"
$ws.Range("A55").ClearFormats()

$ws.Range("B55").Value = "金"
$ws.Range("C55").Value = 9
$ws.Range("D55").Value = 201
